$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet's logical/tab name from "op2" to "wong3"
$ws.Name = "wong3"

# Update row 4 values: B4, C4, E4, F4 from 1 -> 2 (D4 and G4 remain unchanged)
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 2
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 2
